# netCrypto.xlsx - daily update
# Update the USD Amount figure for 11/07/2025 and leave the cursor where
# the analyst left off reviewing the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# T2 ("USD Amount" for the Deposit/Crypto/Roobic row) was revised upward.
$ws.Range("T2").Value = 238727

# Move the active selection to where editing finished (S29).
$ws.Range("S29").Select() | Out-Null

# Best-effort: also nudge the saved window position to match the author's
# last on-screen layout (not all hosts persist window geometry).
try { $excel.Windows.Item(1).Left = -120 } catch {}
